# "Generate Report for Handoff" - refresh the localization-status report with the
# results of a new handoff run (new source-file GUID, new xliff package names,
# refreshed handoff timestamps, and cleared handback info since handback has not
# happened yet for this run).

$wb = $excel.ActiveWorkbook

$oldGuid = "c25b457e-104a-4642-ac03-3b82ba220be7"
$newGuid = "f4f93ec9-d195-4c95-9fae-880eb824291f"
$oldXliffHash = "374c7541cf42aece8acdb28e821c6cb35b10a26e"
$newXliffHash = "24bf2a238efab5f6764214e844f4eea9a9a62130"

$newFileName = "$newGuid.md"
$newPathAndName = "e2e\$newGuid.md"

$zhHandoffFile = "$newGuid.$newXliffHash.zh-cn.xlf"
$deHandoffFile = "$newGuid.$newXliffHash.de-de.xlf"

$zhHandoffDate = "2016-09-06 07:10:52"
$deHandoffDate = "2016-09-06 07:10:58"
$neverHandedBack = "0001-01-01 00:00:00"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b2e5fdcafc0a9de1bbb3781205fa7679c1b0a24/e2e/$newGuid.md"
$zhRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b7a7a9c7daa9866fb2791046928009cddf1b30ac/e2e/$newGuid.md"
$deRepoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/76d5f11a464dac4748a8de06dc082994c2600e73/e2e/$newGuid.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFileName

$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Range("B2").Value = $newPathAndName
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $repoBase, "", "", $newPathAndName) | Out-Null

# Latest HO Xliff Generate Date mirrors the most recent per-language handoff (de-de)
$wsOverview.Range("G2").Value = $deHandoffDate

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = $newFileName
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhRepoBase, "", "", $newFileName) | Out-Null

$wsZh.Range("G2").Value = $zhHandoffFile
$wsZh.Range("H2").Value = $zhHandoffDate

# Latest Target File / Latest Handback File hyperlink+text removed (no handback yet)
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""

# Latest Handback DateTime -> sentinel "never" date
$wsZh.Range("K2").Value = $neverHandedBack

$wsZh.Columns.Item(9).ColumnWidth = 17.8
$wsZh.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = $newFileName
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deRepoBase, "", "", $newFileName) | Out-Null

$wsDe.Range("G2").Value = $deHandoffFile
$wsDe.Range("H2").Value = $deHandoffDate

$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""

$wsDe.Range("K2").Value = $neverHandedBack

$wsDe.Columns.Item(9).ColumnWidth = 17.8
$wsDe.Columns.Item(10).ColumnWidth = 20.8
